$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, matching style of existing header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Add the Save column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
